# Insert a new weekly price record as row 38, pushing the existing
# rows 38-40 down to 39-41 (same as the authored diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 38 (shifts 38..40 down to 39..41)
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new record's data
$ws.Cells.Item(38, 1).Value = 7
$ws.Cells.Item(38, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38, 3).Value = "Ñuble"
$ws.Cells.Item(38, 4).Value = 44610
$ws.Cells.Item(38, 4).NumberFormat = $ws.Cells.Item(39, 4).NumberFormat
$ws.Cells.Item(38, 5).Value = 16
$ws.Cells.Item(38, 6).Value = "Fruta"
$ws.Cells.Item(38, 7).Value = 100103
$ws.Cells.Item(38, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(38, 9).Value = 100103002
$ws.Cells.Item(38, 10).Value = "Ciruela"
$ws.Cells.Item(38, 11).Value = "Black Amber"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 60
$ws.Cells.Item(38, 14).Value = 10000
$ws.Cells.Item(38, 15).Value = 11000
$ws.Cells.Item(38, 16).Value = 10500
$ws.Cells.Item(38, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(38, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(38, 19).Value = 583
$ws.Cells.Item(38, 20).Value = 18
